$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.471.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.946.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.78%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.25%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.939.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.95%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.82"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.87%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.56%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.12%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.75"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.97%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.424.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.96%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +10.24%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.935.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.95%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.456.86"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "417.38"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.44%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.13"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.28%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.16%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.54"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.69%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.30%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.37%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.65%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.102"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.61%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.64"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.18%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.89%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.53"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0687"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.47"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.68%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.92%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.108"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0349"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.20%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "379.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.30%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.652.89"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.02%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.69"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.64%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.79%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.66%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.37"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.56%  "
